# Applies cryptos.xlsx update: refreshed prices/volumes and two row-order swaps (rows 15/16, 29/30)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.957.62'
$ws.Range('E2').Value = '  -5.11%  '
$ws.Range('D3').Value = '2.218.56'
$ws.Range('E3').Value = '  -6.27%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '320.33'
$ws.Range('E5').Value = '  +0.62%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '98.58'
$ws.Range('E6').Value = '  -8.97%  '
$ws.Range('E7').Value = '  -8.66%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.559'
$ws.Range('E9').Value = '  -8.84%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.70'
$ws.Range('E10').Value = '  -10.30%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '53.99'
$ws.Range('E11').Value = '  -3.72%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0825'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '7.61'
$ws.Range('E13').Value = '  -10.37%  '
$ws.Range('E14').Value = '  -1.83%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.863'
$ws.Range('E15').Value = '  -11.99%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '2.553.86'
$ws.Range('E16').Value = '  -6.36%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.32'
$ws.Range('E17').Value = '  -7.04%  '
$ws.Range('D18').Value = '2.219.38'
$ws.Range('E18').Value = '  -6.24%  '
$ws.Range('D19').Value = '42.839.84'
$ws.Range('E19').Value = '  -5.19%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.00'
$ws.Range('E20').Value = '  -9.83%  '
$ws.Range('D21').Value = '0.0₃0963'
$ws.Range('E21').Value = '  -9.46%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.53'
$ws.Range('E22').Value = '  -10.51%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '65.06'
$ws.Range('E23').Value = '  -11.08%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.19'
$ws.Range('E24').Value = '  -11.46%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '236.30'
$ws.Range('E25').Value = '  -10.92%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.15'
$ws.Range('E26').Value = '  -7.81%  '
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('E28').Value = '  +1.08%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.94'
$ws.Range('E29').Value = '  -11.58%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.17'
$ws.Range('E30').Value = '  -5.97%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.33'
$ws.Range('E31').Value = '  -15.54%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '35.58'
$ws.Range('E32').Value = '  -4.73%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.24'
$ws.Range('E33').Value = '  -9.53%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0864'
$ws.Range('E34').Value = '  -8.72%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '153.62'
$ws.Range('E35').Value = '  -9.07%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.65'
$ws.Range('E36').Value = '  -7.92%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.21'
$ws.Range('E37').Value = '  +3.89%  '
$ws.Range('E38').Value = '  -7.92%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.91'
$ws.Range('E39').Value = '  -0.92%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.39'
$ws.Range('E40').Value = '  -6.83%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.104'
$ws.Range('E41').Value = '  -10.91%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.67'
$ws.Range('E42').Value = '  -9.50%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0320'
$ws.Range('E43').Value = '  -9.44%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.70'
$ws.Range('E44').Value = '  +5.55%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = '1.725.83'
$ws.Range('E46').Value = '  -8.28%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.203'
$ws.Range('E47').Value = '  -10.78%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '84.41'
$ws.Range('E48').Value = '  -15.13%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.26'
$ws.Range('E49').Value = '  -11.92%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.77'
$ws.Range('E50').Value = '  -4.44%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '74.23'
$ws.Range('E51').Value = '  -12.04%  '
